$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# The "Enterprises (absolute #)" row and the "Enterprises density (per 1000 people)"
# row (and their data) swap places in both indicator tables on the sheet.
# Use Copy / PasteSpecial(xlPasteValues) so the numeric-looking labels stay stored
# as text (matching their original cell type) and keep their original cell style -
# a plain Range.Value assignment would get auto-coerced to a number and a
# NumberFormat workaround would leave stray style entries behind.

# --- Table 1 (rows 10-11, columns A:D): swap the two entire rows ---
$ws.Range("A10:D10").Copy() | Out-Null
$ws.Range("F1:I1").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A10:D10").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("F1:I1").Copy() | Out-Null
$ws.Range("A11:D11").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("F1:I1").ClearContents() | Out-Null

# --- Table 2 (rows 29-30): swap the label (col A) and data (col D) ---
$ws.Range("A29").Copy() | Out-Null
$ws.Range("F1").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("A30").Copy() | Out-Null
$ws.Range("A29").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("F1").Copy() | Out-Null
$ws.Range("A30").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("D29").Copy() | Out-Null
$ws.Range("F1").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("D30").Copy() | Out-Null
$ws.Range("D29").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("F1").Copy() | Out-Null
$ws.Range("D30").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("F1").ClearContents() | Out-Null
